# Update crypto price (D) and 1h-volume-change (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.574.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.75"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9950"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.89%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "282.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9967"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5137"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3536"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.18"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06872"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.14"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8175"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07763"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.864.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.43"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.123"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9937"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.33"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008135"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9964"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.571.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.814"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.14"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.243"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.400"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +12.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.53"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.665"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.32"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.61"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.389"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.339"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08800"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04913"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.178"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7456"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.867"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.273"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.417"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01870"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5232"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9672"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "116.55"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.307"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.093"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9964"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4579"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1366"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.435"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.52"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.514"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05927"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.83%  "

